# "column full isn't working"
# Fill in the timesheet row for 3/14/2020 (row 15, "Finnished CheckWin")
# estimated-time / remaining-time columns (G/H), and add the missing
# 3/15/2020 entry (row 16) for "Worked on game logic".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Record")

# --- Row 15: fill in the Estimated Time (G) and Remaining (H) columns ---
$ws.Range("G15").Value = 3.5
$ws.Range("H15").Formula = "=G15-D15"

# --- Row 16: brand new log entry ---
$ws.Range("A16").NumberFormat = "d-mmm"
$ws.Range("A16").Value = 43905

$ws.Range("B16").NumberFormat = "h:mm"
$ws.Range("C16").NumberFormat = "h:mm"
$ws.Range("B16").Value = 0.52430555555555558
$ws.Range("C16").Value = 0.59375

$ws.Range("E16").Value = "Worked on game logic"
$ws.Range("F16").Value = "1d"

$ws.Range("G16").Value = 3.5
$ws.Range("H16").Formula = "=H15-D16"

# --- Update the active selection to reflect where the author left off ---
$null = $ws.Range("B17").Select()
